$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.414.16'
$ws.Range('D3').Value = '1.849.86'
$ws.Range('E3').Value = '  +0.12%  '
$r = $ws.Range('D4')
$r.NumberFormat = '@'
$r.Value = '1.000'
$r.ClearFormats()
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('E7').Value = '  +0.09%  '
$r = $ws.Range('D8')
$r.NumberFormat = '@'
$r.Value = '0.07652'
$r.ClearFormats()
$ws.Range('E8').Value = '  +0.56%  '
$r = $ws.Range('D9')
$r.NumberFormat = '@'
$r.Value = '0.2910'
$r.ClearFormats()
$ws.Range('E9').Value = '  -0.74%  '
$ws.Range('E10').Value = '  +1.40%  '
$ws.Range('D11').Value = '2.239.31'
$ws.Range('E11').Value = '  +21.22%  '
$r = $ws.Range('D12')
$r.NumberFormat = '@'
$r.Value = '0.07737'
$r.ClearFormats()
$ws.Range('E12').Value = '  -0.08%  '
$r = $ws.Range('D13')
$r.NumberFormat = '@'
$r.Value = '5.040'
$r.ClearFormats()
$ws.Range('E13').Value = '  +0.73%  '
$r = $ws.Range('D14')
$r.NumberFormat = '@'
$r.Value = '0.6811'
$r.ClearFormats()
$ws.Range('E14').Value = '  +0.24%  '
$r = $ws.Range('D15')
$r.NumberFormat = '@'
$r.Value = '0.00001072'
$r.ClearFormats()
$ws.Range('E15').Value = '  -1.62%  '
$r = $ws.Range('D16')
$r.NumberFormat = '@'
$r.Value = '83.35'
$r.ClearFormats()
$ws.Range('E16').Value = '  -0.50%  '
$r = $ws.Range('D17')
$r.NumberFormat = '@'
$r.Value = '6.173'
$r.ClearFormats()
$ws.Range('E17').Value = '  -0.23%  '
$ws.Range('D18').Value = '29.464.17'
$ws.Range('E18').Value = '  +0.16%  '
$r = $ws.Range('D19')
$r.NumberFormat = '@'
$r.Value = '228.21'
$r.ClearFormats()
$ws.Range('E19').Value = '  -0.27%  '
$ws.Range('E20').Value = '  -0.94%  '
$ws.Range('E21').Value = '  +0.09%  '
$r = $ws.Range('D22')
$r.NumberFormat = '@'
$r.Value = '7.467'
$r.ClearFormats()
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('E23').Value = '  +0.06%  '
$r = $ws.Range('D24')
$r.NumberFormat = '@'
$r.Value = '158.00'
$r.ClearFormats()
$ws.Range('E24').Value = '  +0.38%  '
$r = $ws.Range('D25')
$r.NumberFormat = '@'
$r.Value = '0.1382'
$r.ClearFormats()
$ws.Range('E25').Value = '  -1.16%  '
$r = $ws.Range('D26')
$r.NumberFormat = '@'
$r.Value = '8.421'
$r.ClearFormats()
$ws.Range('E26').Value = '  +0.71%  '
$ws.Range('E27').Value = '  +0.40%  '
$ws.Range('E28').Value = '  +7.05%  '
$r = $ws.Range('D29')
$r.NumberFormat = '@'
$r.Value = '1.462'
$r.ClearFormats()
$ws.Range('E29').Value = '  -0.25%  '
$ws.Range('E30').Value = '  +0.28%  '
$r = $ws.Range('D31')
$r.NumberFormat = '@'
$r.Value = '4.132'
$r.ClearFormats()
$ws.Range('E31').Value = '  +0.74%  '
$r = $ws.Range('D32')
$r.NumberFormat = '@'
$r.Value = '4.065'
$r.ClearFormats()
$ws.Range('E32').Value = '  +0.85%  '
$ws.Range('E34').Value = '  +0.71%  '
$r = $ws.Range('D35')
$r.NumberFormat = '@'
$r.Value = '0.6937'
$r.ClearFormats()
$ws.Range('E35').Value = '  -2.31%  '
$ws.Range('E36').Value = '  +0.18%  '
$r = $ws.Range('D37')
$r.NumberFormat = '@'
$r.Value = '0.01803'
$r.ClearFormats()
$ws.Range('E37').Value = '  +0.18%  '
$ws.Range('D38').Value = '1.229.64'
$ws.Range('E38').Value = '  -0.32%  '
$r = $ws.Range('D39')
$r.NumberFormat = '@'
$r.Value = '2.727'
$r.ClearFormats()
$ws.Range('E39').Value = '  -1.71%  '
$r = $ws.Range('D40')
$r.NumberFormat = '@'
$r.Value = '6.445'
$r.ClearFormats()
$ws.Range('E40').Value = '  +0.18%  '
$r = $ws.Range('D41')
$r.NumberFormat = '@'
$r.Value = '0.9098'
$r.ClearFormats()
$ws.Range('E41').Value = '  +0.23%  '
$ws.Range('E42').Value = '  +0.12%  '
$r = $ws.Range('D43')
$r.NumberFormat = '@'
$r.Value = '101.74'
$r.ClearFormats()
$ws.Range('E43').Value = '  -0.15%  '
$r = $ws.Range('D44')
$r.NumberFormat = '@'
$r.Value = '66.03'
$r.ClearFormats()
$ws.Range('E44').Value = '  -0.14%  '
$r = $ws.Range('D45')
$r.NumberFormat = '@'
$r.Value = '7.198'
$r.ClearFormats()
$ws.Range('E45').Value = '  +0.33%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$r = $ws.Range('D46')
$r.NumberFormat = '@'
$r.Value = '0.00000000118'
$r.ClearFormats()
$ws.Range('E46').Value = '  -2.29%  '
$ws.Range('B47').Value = 'TheSandbox'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$r = $ws.Range('D47')
$r.NumberFormat = '@'
$r.Value = '0.4024'
$r.ClearFormats()
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$r = $ws.Range('D48')
$r.NumberFormat = '@'
$r.Value = '9.035'
$r.ClearFormats()
$ws.Range('E48').Value = '  +0.54%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$r = $ws.Range('D49')
$r.NumberFormat = '@'
$r.Value = '1.684'
$r.ClearFormats()
$ws.Range('E49').Value = '  +0.32%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$r = $ws.Range('D50')
$r.NumberFormat = '@'
$r.Value = '0.1149'
$r.ClearFormats()
$ws.Range('E50').Value = '  +2.32%  '
$r = $ws.Range('D51')
$r.NumberFormat = '@'
$r.Value = '0.05701'
$r.ClearFormats()
$ws.Range('E51').Value = '  -0.14%  '
